# Update "想去人数" (number of people wanting to go) counts on the
# "展览" and "全部类型" sheets to reflect the newly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value to apply on each of the target sheets.
$updates = @{
    "F3"  = 13787
    "F15" = 5750
    "F16" = 127
    "F17" = 85
    "F19" = 83
    "F21" = 148
    "F22" = 224
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
